# Re-upload of database/generate_data/quiz_question.xlsx
#
# The sheet held a "0000NN"-style QUESTION_ID column and QUIZ_ID values in the
# 73000s; this edit replaces them with small sequential integers (QUESTION_ID
# 1..16, QUIZ_ID 1..8) and drops the trailing blank row (row 18).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row([int]$row, [string]$qid, [int]$quizId, [string]$question, [string]$a1, [string]$a2, [string]$a3, [string]$a4, $correct, $deleted) {
    $ws.Cells.Item($row, 1).Value = $qid
    $ws.Cells.Item($row, 2).Value = $quizId
    $ws.Cells.Item($row, 3).Value = $question
    $ws.Cells.Item($row, 4).Value = $a1
    $ws.Cells.Item($row, 5).Value = $a2
    if ($a3 -ne $null) { $ws.Cells.Item($row, 6).Value = $a3 }
    if ($a4 -ne $null) { $ws.Cells.Item($row, 7).Value = $a4 }
    if ($correct -ne $null) { $ws.Cells.Item($row, 8).Value = $correct }
    if ($deleted -ne $null) { $ws.Cells.Item($row, 9).Value = $deleted }
}

# Quiz 1 (QUIZ_ID 1): attendance + "print Hello world" in Python
Set-Row 2 "1" 1 "Bạn có đi học ngày hôm nay không?" "Có" "Không" $null $null $null $null
Set-Row 3 "2" 1 'Viết câu lệnh in ra màn hình "Hello world!" bằng Python.' 'sys.out("Hello World")' 'print("Hello World")' 'cout << "Hello World!";' 'cin >> "Hello world!";' 1 0

# Quiz 2 (QUIZ_ID 2): attendance + "print Hello world" in Python
Set-Row 4 "3" 2 "Bạn có đi học ngày hôm nay không?" "Có" "Không" $null $null $null $null
Set-Row 5 "4" 2 'Viết câu lệnh in ra màn hình "Hello world!" bằng Python.' 'sys.out("Hello World")' 'print("Hello World")' 'cout << "Hello World!";' 'cin >> "Hello world!";' 1 0

# Quiz 3 (QUIZ_ID 3): attendance + sum of two variables a, b
Set-Row 6 "5" 3 "Bạn có đi học ngày hôm nay không?" "Có" "Không" $null $null $null $null
Set-Row 7 "6" 3 "Câu lệnh xuất tổng của 2 biến a, b trong Python là:" "cout << a+b;" "print(a+b)" "cin >> a+b;" "sys.console.print(a+b);" 1 0

# Quiz 4 (QUIZ_ID 4): attendance + for-loop running n times
Set-Row 8 "7" 4 "Bạn có đi học ngày hôm nay không?" "Có" "Không" $null $null $null $null
Set-Row 9 "8" 4 "Câu lệnh cho vòng lặp for chạy n lần trong Python là:" "for i in range(n):" "for i in range(1,n):" "for _ in range(1,n)" "for hehe in list(n):" 0 0

# Quiz 5 (QUIZ_ID 5): attendance + select all rows from mysql "user" table
Set-Row 10 "9"  5 "Bạn có đi học ngày hôm nay không?" "Có" "Không" $null $null $null $null
Set-Row 11 "10" 5 "Câu lệnh chọn tất cả các dòng trong bảng user của mysql là:" "SELECT ALL FROM USER" "SELECT ALL ABOUT USER" "SELECT * FROM USER" "SELECT X FROM USER" 2 0

# Quiz 6 (QUIZ_ID 6): attendance + select all rows from mysql "user" table
Set-Row 12 "11" 6 "Bạn có đi học ngày hôm nay không?" "Có" "Không" $null $null $null $null
Set-Row 13 "12" 6 "Câu lệnh chọn tất cả các dòng trong bảng user của mysql là:" "SELECT ALL FROM USER" "SELECT ALL ABOUT USER" "SELECT * FROM USER" "SELECT X FROM USER" 2 0

# Quiz 7 (QUIZ_ID 7): attendance + "how is your day today"
Set-Row 14 "13" 7 "Bạn có đi học ngày hôm nay không?" "Có" "Không" $null $null $null $null
Set-Row 15 "14" 7 "Ngày hôm nay của bạn thế nào?" "Bình thường" "Bí mật" "Tốt" "Tệ" 1 0

# Quiz 8 (QUIZ_ID 8): attendance + course evaluation
Set-Row 16 "15" 8 "Bạn có đi học ngày hôm nay không?" "Có" "Không" $null $null $null $null
Set-Row 17 "16" 8 "Bạn đánh giá môn học này như thế nào?" "Quá tệ" "Tệ" "Cũng được" "Tốt" 3 0

# The old sheet had a trailing fully-blank row 18; drop it.
$ws.Rows("18").Delete()

# Match the saved selection recorded in the workbook.
$ws.Range("B18").Select()
